$d = $word.ActiveDocument

$pairs = @(
    @("867×8=6936", "543×8=4344"),
    @("703×7=4921", "670×2=1340"),
    @("749×6=4494", "651×3=1953"),
    @("432×6=2592", "444×5=2220"),
    @("111×7=777",  "454×6=2724"),
    @("401×4=1604", "673×3=2019"),
    @("209×6=1254", "494×7=3458"),
    @("643×5=3215", "670×2=1340"),
    @("661×7=4627", "719×5=3595"),
    @("150×2=300",  "157×3=471"),
    @("439×8=3512", "930×8=7440"),
    @("184×2=368",  "602×7=4214"),
    @("443×2=886",  "696×3=2088"),
    @("309×5=1545", "913×8=7304"),
    @("805×4=3220", "755×9=6795"),
    @("905×6=5430", "916×7=6412"),
    @("549×7=3843", "849×4=3396"),
    @("476×3=1428", "567×4=2268"),
    @("435×4=1740", "547×8=4376"),
    @("255×3=765",  "120×5=600"),
    @("890×4=3560", "365×9=3285"),
    @("872×7=6104", "335×4=1340"),
    @("804×9=7236", "993×3=2979"),
    @("574×9=5166", "451×5=2255"),
    @("322×4=1288", "354×2=708")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $new, 2)
}
